# "add special final round zoom"
# 1) Bump the cached DateAndTime placeholder text (1/4/2023 -> 1/6/2023)
#    everywhere it appears: on every slide master and every one of its
#    custom (slide) layouts.
# 2) Update the "Miami University" raffle-room callout on slide 2 to
#    reflect the newly added final-round rooms (10, 11 -> 19, 20).

$p = $ppt.ActivePresentation

$oldDate = "1/4/2023"
$newDate = "1/6/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shape) {
    if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $design = $p.Designs.Item($di)
    $master = $design.SlideMaster

    # The date placeholder living directly on the slide master.
    for ($si = 1; $si -le $master.Shapes.Count; $si++) {
        Update-DatePlaceholder $master.Shapes.Item($si)
    }

    # Each custom layout under this master keeps its own cached copy.
    for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
        $layout = $master.CustomLayouts.Item($li)
        for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
            Update-DatePlaceholder $layout.Shapes.Item($si)
        }
    }
}

# Slide 2 ("Map"): Raffle Rooms list - Miami University now spans two
# additional raffle rooms for the special final-round zoom.
$slide = $p.Slides.Item(2)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $shape = $slide.Shapes.Item($si)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text
        $oldLine = "Miami University: 10, 11"
        $newLine = "Miami University: 19, 20"
        $idx = $fullText.IndexOf($oldLine)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $oldLine.Length)
            $sub.Text = $newLine
        }
    }
}
